$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# Remove the duplicate icon groups: the original "Group 7" and its three
# "Group 7_Copy" duplicates (ids 8, 123, 131, 139 in the source OOXML).
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Group 7" -or $sh.Name -eq "Group 7_Copy") {
        $sh.Delete()
    }
}

# Reposition the remaining "Group 147" icon to the spot previously occupied
# by one of the removed duplicates.
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Group 147") {
        $sh.Left = 159.7336220472441
        $sh.Top = 199.64740757480317
    }
}
